$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.441.58'
$ws.Range('E2').Value = '  +7.84%  '
$ws.Range('D3').Value = '3.402.54'
$ws.Range('E3').Value = '  +4.63%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '412.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '124.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +15.07%  '
$ws.Range('D7').Value = '3.397.75'
$ws.Range('E7').Value = '  +4.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.584'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.640'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.124'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +30.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '41.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.11%  '
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = '3.943.65'
$ws.Range('E14').Value = '  +4.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('D17').Value = '3.403.17'
$ws.Range('E17').Value = '  +4.80%  '
$ws.Range('D18').Value = '61.471.80'
$ws.Range('E18').Value = '  +8.24%  '
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000122'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +12.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '299.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '76.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.83%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '30.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.12%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.117'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.27%  '
$ws.Range('E31').Value = '  +1.87%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.31%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +19.72%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.38'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0477'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.70%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.51'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('E43').Value = '  +1.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '134.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.281'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.18'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '21.77'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.68%  '
$ws.Range('D50').Value = '2.197.33'
$ws.Range('E50').Value = '  +2.45%  '
$ws.Range('D51').Value = '3.745.91'
$ws.Range('E51').Value = '  +4.96%  '
